$d = $word.ActiveDocument
$d.Content.Find.Execute("Its ultimate purpose is the complete mastery of mind over the material world, the harnessing of the forces of nature to human needs. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
